$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Correct the Zone column (B4) for agent "Brann" (row 4): it mistakenly held
# the agent's name again instead of a zone code - set it to "PRM2".
$ws.Range("B4").Value = "PRM2"

# Move the active selection to B5, as it ended up after the edit.
$ws.Range("B5").Select()
